$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap rows 50 and 51 (B:AC) ---
$ws.Range("B50").Value = 6732794
$ws.Range("C50").Value = "Lithuania A Lyga"
$ws.Range("D50").Value = "Lithuania A Lyga"
$ws.Range("E50").Value = 45149.54166666666
$ws.Range("F50").Value = "FK Siauliai"
$ws.Range("G50").Value = "FK Dziugas Telsiai"
$ws.Range("H50").Value = 3
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = "H"
$ws.Range("K50").Value = 1.25
$ws.Range("L50").Value = 5
$ws.Range("M50").Value = 9
$ws.Range("N50").Value = 1.25
$ws.Range("O50").Value = 5.25
$ws.Range("P50").Value = 9
$ws.Range("Q50").Value = -1.75
$ws.Range("R50").Value = 2
$ws.Range("S50").Value = 1.8
$ws.Range("T50").Value = 3
$ws.Range("U50").Value = 1.975
$ws.Range("V50").Value = 1.825
$ws.Range("W50").Value = 0.25
$ws.Range("X50").Value = -1
$ws.Range("Y50").Value = -1
$ws.Range("Z50").Value = 1
$ws.Range("AA50").Value = -1
$ws.Range("AB50").Value = 0
$ws.Range("AC50").Value = 0
$ws.Range("B51").Value = 6732795
$ws.Range("C51").Value = "Lithuania A Lyga"
$ws.Range("D51").Value = "Lithuania A Lyga"
$ws.Range("E51").Value = 45149.54166666666
$ws.Range("F51").Value = "Suduva Marijampole"
$ws.Range("G51").Value = "Banga Gargzdai"
$ws.Range("H51").Value = 1
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = "H"
$ws.Range("K51").Value = 2.15
$ws.Range("L51").Value = 3.2
$ws.Range("M51").Value = 3
$ws.Range("N51").Value = 2.3
$ws.Range("O51").Value = 3.2
$ws.Range("P51").Value = 2.7
$ws.Range("Q51").Value = -0.25
$ws.Range("R51").Value = 2.05
$ws.Range("S51").Value = 1.75
$ws.Range("T51").Value = 2.25
$ws.Range("U51").Value = 1.9
$ws.Range("V51").Value = 1.9
$ws.Range("W51").Value = 1.3
$ws.Range("X51").Value = -1
$ws.Range("Y51").Value = -1
$ws.Range("Z51").Value = 1.05
$ws.Range("AA51").Value = -1
$ws.Range("AB51").Value = -1
$ws.Range("AC51").Value = 0.8999999999999999

# --- Rotate rows 101, 102, 104 (B:AC): 104->101, 101->102, 102->104 ---
$ws.Range("B101").Value = 6732834
$ws.Range("C101").Value = "Lithuania A Lyga"
$ws.Range("D101").Value = "Lithuania A Lyga"
$ws.Range("E101").Value = 45242.41319444445
$ws.Range("F101").Value = "Panevezys"
$ws.Range("G101").Value = "FK Dziugas Telsiai"
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = "D"
$ws.Range("K101").Value = 1.25
$ws.Range("L101").Value = 5.5
$ws.Range("M101").Value = 7.5
$ws.Range("N101").Value = 1.45
$ws.Range("O101").Value = 4.5
$ws.Range("P101").Value = 5
$ws.Range("Q101").Value = -1
$ws.Range("R101").Value = 1.775
$ws.Range("S101").Value = 2.025
$ws.Range("T101").Value = 2.5
$ws.Range("U101").Value = 1.875
$ws.Range("V101").Value = 1.925
$ws.Range("W101").Value = -1
$ws.Range("X101").Value = 3.5
$ws.Range("Y101").Value = -1
$ws.Range("Z101").Value = -1
$ws.Range("AA101").Value = 1.025
$ws.Range("AB101").Value = -1
$ws.Range("AC101").Value = 0.925
$ws.Range("B102").Value = 7465686
$ws.Range("C102").Value = "Lithuania A Lyga"
$ws.Range("D102").Value = "Lithuania A Lyga"
$ws.Range("E102").Value = 45242.41319444445
$ws.Range("F102").Value = "FK Kauno Zalgiris"
$ws.Range("G102").Value = "Hegelmann Litauen"
$ws.Range("H102").Value = 4
$ws.Range("I102").Value = 2
$ws.Range("J102").Value = "H"
$ws.Range("K102").Value = 2.3
$ws.Range("L102").Value = 4
$ws.Range("M102").Value = 2.3
$ws.Range("N102").Value = 2.55
$ws.Range("O102").Value = 4
$ws.Range("P102").Value = 2.2
$ws.Range("Q102").Value = 0.25
$ws.Range("R102").Value = 1.8
$ws.Range("S102").Value = 2
$ws.Range("T102").Value = 2.75
$ws.Range("U102").Value = 1.85
$ws.Range("V102").Value = 1.95
$ws.Range("W102").Value = 1.55
$ws.Range("X102").Value = -1
$ws.Range("Y102").Value = -1
$ws.Range("Z102").Value = 0.8
$ws.Range("AA102").Value = -1
$ws.Range("AB102").Value = 0.8500000000000001
$ws.Range("AC102").Value = -1
$ws.Range("B104").Value = 6732837
$ws.Range("C104").Value = "Lithuania A Lyga"
$ws.Range("D104").Value = "Lithuania A Lyga"
$ws.Range("E104").Value = 45242.41319444445
$ws.Range("F104").Value = "Suduva Marijampole"
$ws.Range("G104").Value = "FK Riteriai"
$ws.Range("H104").Value = 0
$ws.Range("I104").Value = 3
$ws.Range("J104").Value = "A"
$ws.Range("K104").Value = 3.6
$ws.Range("L104").Value = 3.6
$ws.Range("M104").Value = 1.8
$ws.Range("N104").Value = 3
$ws.Range("O104").Value = 3.6
$ws.Range("P104").Value = 2
$ws.Range("Q104").Value = 0.25
$ws.Range("R104").Value = 2
$ws.Range("S104").Value = 1.8
$ws.Range("T104").Value = 2.5
$ws.Range("U104").Value = 1.975
$ws.Range("V104").Value = 1.825
$ws.Range("W104").Value = -1
$ws.Range("X104").Value = -1
$ws.Range("Y104").Value = 1
$ws.Range("Z104").Value = -1
$ws.Range("AA104").Value = 0.8
$ws.Range("AB104").Value = 0.9750000000000001
$ws.Range("AC104").Value = -1

# --- Copy cell formatting from row 130 onto new rows 131:134 ---
$ws.Range("A130:AC130").Copy() | Out-Null
$ws.Range("A131:AC134").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Row 130: replace all values (B:AC), including new H/I/J/AB/AC cells ---
$ws.Range("B130").Value = 7862918
$ws.Range("C130").Value = "Lithuania A Lyga"
$ws.Range("D130").Value = "Lithuania A Lyga"
$ws.Range("E130").Value = 45388.375
$ws.Range("F130").Value = "FK Dziugas Telsiai"
$ws.Range("G130").Value = "Hegelmann Litauen"
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = "D"
$ws.Range("K130").Value = 4.333
$ws.Range("L130").Value = 3.75
$ws.Range("M130").Value = 1.615
$ws.Range("N130").Value = 4.5
$ws.Range("O130").Value = 3.4
$ws.Range("P130").Value = 1.65
$ws.Range("Q130").Value = 0.75
$ws.Range("R130").Value = 1.875
$ws.Range("S130").Value = 1.925
$ws.Range("T130").Value = 2.25
$ws.Range("U130").Value = 1.875
$ws.Range("V130").Value = 1.925
$ws.Range("W130").Value = -1
$ws.Range("X130").Value = 2.4
$ws.Range("Y130").Value = -1
$ws.Range("Z130").Value = 0.875
$ws.Range("AA130").Value = -1
$ws.Range("AB130").Value = -1
$ws.Range("AC130").Value = 0.925

# --- New row 131 ---
$ws.Range("A131").Value = 129
$ws.Range("B131").Value = 7862919
$ws.Range("C131").Value = "Lithuania A Lyga"
$ws.Range("D131").Value = "Lithuania A Lyga"
$ws.Range("E131").Value = 45388.45833333334
$ws.Range("F131").Value = "FK Dainava Alytus"
$ws.Range("G131").Value = "Banga Gargzdai"
$ws.Range("H131").Value = 2
$ws.Range("I131").Value = 2
$ws.Range("J131").Value = "D"
$ws.Range("K131").Value = 2.25
$ws.Range("L131").Value = 2.9
$ws.Range("M131").Value = 3.1
$ws.Range("N131").Value = 2.05
$ws.Range("O131").Value = 3
$ws.Range("P131").Value = 3.5
$ws.Range("Q131").Value = -0.25
$ws.Range("R131").Value = 1.775
$ws.Range("S131").Value = 2.025
$ws.Range("T131").Value = 2
$ws.Range("U131").Value = 1.9
$ws.Range("V131").Value = 1.9
$ws.Range("W131").Value = -1
$ws.Range("X131").Value = 2
$ws.Range("Y131").Value = -1
$ws.Range("Z131").Value = -0.5
$ws.Range("AA131").Value = 0.5125
$ws.Range("AB131").Value = 0.8999999999999999
$ws.Range("AC131").Value = -1

# --- New row 132 ---
$ws.Range("A132").Value = 130
$ws.Range("B132").Value = 7865009
$ws.Range("C132").Value = "Lithuania A Lyga"
$ws.Range("D132").Value = "Lithuania A Lyga"
$ws.Range("E132").Value = 45389.29166666666
$ws.Range("F132").Value = "FK Transinvest"
$ws.Range("G132").Value = "Suduva Marijampole"
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 2
$ws.Range("J132").Value = "A"
$ws.Range("K132").Value = 1.833
$ws.Range("L132").Value = 3.4
$ws.Range("M132").Value = 3.6
$ws.Range("N132").Value = 2
$ws.Range("O132").Value = 3.25
$ws.Range("P132").Value = 3.3
$ws.Range("Q132").Value = -0.25
$ws.Range("R132").Value = 1.775
$ws.Range("S132").Value = 2.025
$ws.Range("T132").Value = 2.25
$ws.Range("U132").Value = 1.875
$ws.Range("V132").Value = 1.925
$ws.Range("W132").Value = -1
$ws.Range("X132").Value = -1
$ws.Range("Y132").Value = 2.3
$ws.Range("Z132").Value = -1
$ws.Range("AA132").Value = 1.025
$ws.Range("AB132").Value = -0.5
$ws.Range("AC132").Value = 0.4625

# --- New row 133 ---
$ws.Range("A133").Value = 131
$ws.Range("B133").Value = 7862043
$ws.Range("C133").Value = "Lithuania A Lyga"
$ws.Range("D133").Value = "Lithuania A Lyga"
$ws.Range("E133").Value = 45389.375
$ws.Range("F133").Value = "FK Zalgiris Vilnius"
$ws.Range("G133").Value = "FK Siauliai"
$ws.Range("H133").Value = 3
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = "H"
$ws.Range("K133").Value = 1.666
$ws.Range("L133").Value = 3.75
$ws.Range("M133").Value = 4
$ws.Range("N133").Value = 1.7
$ws.Range("O133").Value = 3.6
$ws.Range("P133").Value = 4
$ws.Range("Q133").Value = -0.75
$ws.Range("R133").Value = 1.975
$ws.Range("S133").Value = 1.825
$ws.Range("T133").Value = 2.5
$ws.Range("U133").Value = 1.925
$ws.Range("V133").Value = 1.875
$ws.Range("W133").Value = 0.7
$ws.Range("X133").Value = -1
$ws.Range("Y133").Value = -1
$ws.Range("Z133").Value = 0.9750000000000001
$ws.Range("AA133").Value = -1
$ws.Range("AB133").Value = 0.925
$ws.Range("AC133").Value = -1

# --- New row 134 ---
$ws.Range("A134").Value = 132
$ws.Range("B134").Value = 7862920
$ws.Range("C134").Value = "Lithuania A Lyga"
$ws.Range("D134").Value = "Lithuania A Lyga"
$ws.Range("E134").Value = 45389.51736111111
$ws.Range("F134").Value = "FK Kauno Zalgiris"
$ws.Range("G134").Value = "Panevezys"
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 1
$ws.Range("J134").Value = "A"
$ws.Range("K134").Value = 2.75
$ws.Range("L134").Value = 3
$ws.Range("M134").Value = 2.4
$ws.Range("N134").Value = 2.375
$ws.Range("O134").Value = 3
$ws.Range("P134").Value = 2.8
$ws.Range("Q134").Value = 0
$ws.Range("R134").Value = 1.8
$ws.Range("S134").Value = 2
$ws.Range("T134").Value = 2
$ws.Range("U134").Value = 1.85
$ws.Range("V134").Value = 1.95
$ws.Range("W134").Value = -1
$ws.Range("X134").Value = -1
$ws.Range("Y134").Value = 1.8
$ws.Range("Z134").Value = -1
$ws.Range("AA134").Value = 1
$ws.Range("AB134").Value = -1
$ws.Range("AC134").Value = 0.95
